$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "sequence/run_2288_samples/" path from each fastqFileName
# value in column F (rows 2-25), leaving just the bare filename.
$prefix = "sequence/run_2288_samples/"

for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().StartsWith($prefix)) {
        $cell.Value = $val.ToString().Substring($prefix.Length)
    }
}

# Update the active cell/selection on the sheet to match the target state.
$ws.Range("P23").Select()
